# Adds three new literature-review paragraphs (with justified alignment)
# before the existing trailing paragraph, matching the target diff.
$d = $word.ActiveDocument

$newParagraphsXml = @'
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Learning structural objects from unknown probability distribution is becoming popular in recent years. </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>Tsochantaridis</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> et al.</w:t>
  </w:r>
  <w:r>
    <w:t>\cite</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>generalized</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> multiclass SVMs</w:t>
  </w:r>
  <w:r>
    <w:t>\cite</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> to structural SVMs by </w:t>
  </w:r>
  <w:r>
    <w:t>extending feature vectors</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">to joint feature vectors which map features extracted jointly over input-output pairs to </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">discrete output. The </w:t>
  </w:r>
  <w:r>
    <w:t>exact maximum a p</w:t>
  </w:r>
  <w:r>
    <w:t>osteriori</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> (MAP) problem thus becomes an NP-hard problem. They overcome this by using a method called “Soft-Margin Maximization” and found an upper bound of the loss function.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Based on the previous research, Yu and </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>Joachims</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>\cite</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">developed latent SVM by introducing a hidden variable into the joint feature vector. </w:t>
  </w:r>
  <w:r>
    <w:t>They observed a</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>fa</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">ct that in real world applications hidden variables are usually intermediate results and are not required as an output. With this </w:t>
  </w:r>
  <w:r>
    <w:t>fa</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">ct they followed “Soft-Margin” method and found an upper bound for the loss function with latent variables. </w:t>
  </w:r>
  <w:r>
    <w:t>However, the resulted object function is still non-convex.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>Yuille</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> and </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>Rangarajan</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> \cite developed the Concave-Convex Procedure (CCCP) which is guaranteed to find a local minimum for a Difference-Convex (DC) program. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">Yu and </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>Joachims</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve">\cite </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">combined CCCP algorithm by writing their non-convex object function into a difference of two convex function and came up with a EM like 2 steps optimizing algorithm. For each iteration, they first compute latent variables utilizing current parameter vectors and then in turn optimizing parameter vectors using the standard Structural SVM algorithm with previously computed latent variables. </w:t>
  </w:r>
</w:p>
<w:p/>
'@

$insertionPoint = $d.Paragraphs(1).Range
$insertionPoint.Collapse(1)
$insertionPoint.InsertXML($newParagraphsXml)

# InsertXML merges the final inserted paragraph mark into the paragraph that
# followed the insertion point (standard Word paste semantics), so a throwaway
# empty paragraph was appended above to absorb that merge; remove it now that
# the three real paragraphs have their own paragraph marks.
$d.Paragraphs(4).Range.Delete()
